$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix typo in G1 header: "鞋类类..." -> "鞋类..." (duplicated "类" removed)
$ws.Range("G1").Value = "鞋类城市居民消费价格指数(上年=100)"

# Append 2021 data row (row 7), reusing row 6's header-cell style for A7
$ws.Range("A6").Copy($ws.Range("A7"))
$ws.Range("A7").Value = "2021年"
$ws.Range("D7").Value = 100.5
$ws.Range("F7").Value = 100.3
$ws.Range("G7").Value = 99.90000000000001

# Append 2022 data row (row 8), reusing the same style for A8
$ws.Range("A6").Copy($ws.Range("A8"))
$ws.Range("A8").Value = "2022年"
$ws.Range("F8").Value = 100.6
